# Update the two JIRA URL values in column A (shared-string backed cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "https://jira.com/browse/qw123"
$ws.Range("A3").Value = "https://jira.com/browse/qw234"

# Widen column A to fit the new (shorter) text and move the active selection
# down to A4, matching the saved view state.
$ws.Columns("A").ColumnWidth = 40.74

$ws.Range("A4").Select()
